# #5: cash & deposit done
# Rebuild the "存款" (deposit) sheet so that it carries the same full set of
# columns already used on the other property-category sheets (土地/建物/汽車):
#   A index | B bank | C deposit_type | D currency | E owner | F total |
#   G property_category | H category | I date | J legislator_name |
#   K legislator_id | L source_file | M index

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# A sheet that already has a literal "2012-04-30" text value we can clone from,
# so that assigning the date string to column I doesn't get auto-coerced into
# a date serial number by the COM layer.
$wsCar = $wb.Worksheets.Item("汽車")

# ---- Row 1: turn the old (duplicate-of-row-2) data row into a proper header
#      row of field-name labels. ----
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"
$ws.Cells.Item(1,7).Value = "property_category"
$ws.Cells.Item(1,8).Value = "category"
$ws.Cells.Item(1,9).Value = "date"
$ws.Cells.Item(1,10).Value = "legislator_name"
$ws.Cells.Item(1,11).Value = "legislator_id"
$ws.Cells.Item(1,12).Value = "source_file"
$ws.Cells.Item(1,13).Value = "index"

# New header cells G1:M1 should carry the same bold/bordered/centered style as
# the existing header cells B1:F1 -- clone the format in one bulk operation so
# we reuse the existing style record instead of building up new partial ones.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)   # xlPasteFormats

# ---- Data rows 2-8: columns A-F already hold the correct values; fix the one
#      cell that was mistakenly stored as text, then append the new columns
#      G-M that carry property_category/category/date/legislator info. ----

# Row 6 (index 50) column F held the amount as the text "1693963" -- correct
# it to a real number while keeping the existing formatting.
$ws.Cells.Item(6,6).Value = 1693963

$indices = @{2=46; 3=47; 4=48; 5=49; 6=50; 7=51; 8=52}

# Copy a cell that already stores the literal text "2012-04-30" once, then
# paste-values it into every row's date column.
$wsCar.Range("J2").Copy()

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r,7).Value  = "deposit"
    $ws.Cells.Item($r,8).Value  = "normal"
    $ws.Cells.Item($r,9).PasteSpecial(-4163)   # xlPasteValues -> keeps it text "2012-04-30"
    $ws.Cells.Item($r,10).Value = "趙天麟"
    $ws.Cells.Item($r,11).Value = 1761
    $ws.Cells.Item($r,12).Value = "tmp58581"
    $ws.Cells.Item($r,13).Value = $indices[$r]
}

Write-Output "sheet4 (存款) rebuilt with full column set"
